# Jeux.xlsx - "Changement excel liés au trigger + date jeux"
#
# 1) Column F (date_parution) currently stores a bare 4-digit year as an
#    integer (e.g. 2011). Convert every value to a real Excel date
#    (1 January of that year) and display it with a custom "yyyy" number
#    format, so it keeps showing just the year but is stored as a date
#    serial (needed for a date-based trigger/computation).
# 2) A new (currently empty) column J is introduced next to H, formatted
#    with an explicit General number format - this is the column used by
#    the new trigger.
# 3) Minor sheet/view housekeeping: selection, page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 51
$lastCol = 6   # column F

# --- 1) Rewrite column F: year (e.g. 2011) -> date serial for Jan 1st of that year
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $lastCol)
    $year = [int]$cell.Value2
    $jan1 = Get-Date -Year $year -Month 1 -Day 1
    $serial = [Math]::Floor($jan1.ToOADate())
    $cell.Value = $serial
}

# --- 2) New column J: empty cells with an explicit "General" number format
#        (create this style first so it lands at cellXfs index 1, matching
#        the order in which the F-column date format is created next)
$ws.Range("J2:J51").NumberFormat = "General"

# --- Apply the custom "yyyy" display format to the whole F column of data
$ws.Range("F2:F51").NumberFormat = "yyyy"

# --- 3) View / selection tweaks
$ws.Range("I11").Select()

# --- Page setup
$ps = $ws.PageSetup
$ps.PaperSize = 9        # xlPaperA4
$ps.Orientation = 1      # xlPortrait
